$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Hgf"
$ws.Range("C2").Value = "Sdc2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.09364566666666667
$ws.Range("H2").Value = 0.280937
$ws.Range("I2").Value = 0.00132949638239575
$ws.Range("J2").Value = 0.001329496382395751
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.248835333333334
$ws.Range("N2").Value = 6.746506
$ws.Range("O2").Value = 0.03590294220158827
$ws.Range("P2").Value = 0.03590294220158827
$ws.Range("Q2").Value = 0.2105936840135556
$ws.Range("R2").Value = 1.895343156122
$ws.Range("S2").Value = 0.00004773283177437532
$ws.Range("T2").Value = 0.00004773283177437534

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Hgf"
$ws.Range("C3").Value = "Sdc2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.09364566666666667
$ws.Range("H3").Value = 0.280937
$ws.Range("I3").Value = 0.00132949638239575
$ws.Range("J3").Value = 0.001329496382395751
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 44.29005966666667
$ws.Range("N3").Value = 132.870179
$ws.Range("O3").Value = 0.7070964373190639
$ws.Range("P3").Value = 0.7070964373190639
$ws.Range("Q3").Value = 4.147572164191445
$ws.Range("R3").Value = 37.328149477723
$ws.Range("S3").Value = 0.0009400821554206189
$ws.Range("T3").Value = 0.000940082155420619

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Hgf"
$ws.Range("C4").Value = "Sdc2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.09364566666666667
$ws.Range("H4").Value = 0.280937
$ws.Range("I4").Value = 0.00132949638239575
$ws.Range("J4").Value = 0.001329496382395751
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 16.09762433333333
$ws.Range("N4").Value = 48.292873
$ws.Range("O4").Value = 0.2570006204793478
$ws.Range("P4").Value = 0.2570006204793479
$ws.Range("Q4").Value = 1.507472762444555
$ws.Range("R4").Value = 13.567254862001
$ws.Range("S4").Value = 0.0003416813952007561
$ws.Range("T4").Value = 0.0003416813952007562

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Hgf"
$ws.Range("C5").Value = "Sdc2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 10.26742066666667
$ws.Range("H5").Value = 30.802262
$ws.Range("I5").Value = 0.1457675418282607
$ws.Range("J5").Value = 0.1457675418282608
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.248835333333334
$ws.Range("N5").Value = 6.746506
$ws.Range("O5").Value = 0.03590294220158827
$ws.Range("P5").Value = 0.03590294220158827
$ws.Range("Q5").Value = 23.08973837739689
$ws.Range("R5").Value = 207.807645396572
$ws.Range("S5").Value = 0.005233483629127646
$ws.Range("T5").Value = 0.005233483629127647

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Hgf"
$ws.Range("C6").Value = "Sdc2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 10.26742066666667
$ws.Range("H6").Value = 30.802262
$ws.Range("I6").Value = 0.1457675418282607
$ws.Range("J6").Value = 0.1457675418282608
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 44.29005966666667
$ws.Range("N6").Value = 132.870179
$ws.Range("O6").Value = 0.7070964373190639
$ws.Range("P6").Value = 0.7070964373190639
$ws.Range("Q6").Value = 454.7446739494331
$ws.Range("R6").Value = 4092.702065544898
$ws.Range("S6").Value = 0.1030717095035208
$ws.Range("T6").Value = 0.1030717095035208

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Hgf"
$ws.Range("C7").Value = "Sdc2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 10.26742066666667
$ws.Range("H7").Value = 30.802262
$ws.Range("I7").Value = 0.1457675418282607
$ws.Range("J7").Value = 0.1457675418282608
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 16.09762433333333
$ws.Range("N7").Value = 48.292873
$ws.Range("O7").Value = 0.2570006204793478
$ws.Range("P7").Value = 0.2570006204793479
$ws.Range("Q7").Value = 165.2810807643029
$ws.Range("R7").Value = 1487.529726878726
$ws.Range("S7").Value = 0.0374623486956123
$ws.Range("T7").Value = 0.03746234869561231

$ws.Range("A8").Value = "M1"
$ws.Range("B8").Value = "Hgf"
$ws.Range("C8").Value = "Sdc2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 30.40955033333333
$ws.Range("H8").Value = 91.228651
$ws.Range("I8").Value = 0.4317272608283866
$ws.Range("J8").Value = 0.4317272608283867
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 2.248835333333334
$ws.Range("N8").Value = 6.746506
$ws.Range("O8").Value = 0.03590294220158827
$ws.Range("P8").Value = 0.03590294220158827
$ws.Range("Q8").Value = 68.38607126037844
$ws.Range("R8").Value = 615.474641343406
$ws.Range("S8").Value = 0.01550027889237159
$ws.Range("T8").Value = 0.01550027889237159

$ws.Range("A9").Value = "M1"
$ws.Range("B9").Value = "Hgf"
$ws.Range("C9").Value = "Sdc2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 30.40955033333333
$ws.Range("H9").Value = 91.228651
$ws.Range("I9").Value = 0.4317272608283866
$ws.Range("J9").Value = 0.4317272608283867
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 44.29005966666667
$ws.Range("N9").Value = 132.870179
$ws.Range("O9").Value = 0.7070964373190639
$ws.Range("P9").Value = 0.7070964373190639
$ws.Range("Q9").Value = 1346.840798699837
$ws.Range("R9").Value = 12121.56718829853
$ws.Range("S9").Value = 0.3052728080252704
$ws.Range("T9").Value = 0.3052728080252705

$ws.Range("A10").Value = "M1"
$ws.Range("B10").Value = "Hgf"
$ws.Range("C10").Value = "Sdc2"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 30.40955033333333
$ws.Range("H10").Value = 91.228651
$ws.Range("I10").Value = 0.4317272608283866
$ws.Range("J10").Value = 0.4317272608283867
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 16.09762433333333
$ws.Range("N10").Value = 48.292873
$ws.Range("O10").Value = 0.2570006204793478
$ws.Range("P10").Value = 0.2570006204793479
$ws.Range("Q10").Value = 489.5215174115914
$ws.Range("R10").Value = 4405.693656704323
$ws.Range("S10").Value = 0.1109541739107446
$ws.Range("T10").Value = 0.1109541739107446

$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Hgf"
$ws.Range("C11").Value = "Sdc2"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 29.37284333333333
$ws.Range("H11").Value = 88.11852999999999
$ws.Range("I11").Value = 0.417009033544999
$ws.Range("J11").Value = 0.4170090335449991
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 2.248835333333334
$ws.Range("N11").Value = 6.746506
$ws.Range("O11").Value = 0.03590294220158827
$ws.Range("P11").Value = 0.03590294220158827
$ws.Range("Q11").Value = 66.05468792846445
$ws.Range("R11").Value = 594.49219135618
$ws.Range("S11").Value = 0.01497185122890629
$ws.Range("T11").Value = 0.01497185122890629

$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Hgf"
$ws.Range("C12").Value = "Sdc2"
$ws.Range("D12").Value = "FAPs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 29.37284333333333
$ws.Range("H12").Value = 88.11852999999999
$ws.Range("I12").Value = 0.417009033544999
$ws.Range("J12").Value = 0.4170090335449991
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 44.29005966666667
$ws.Range("N12").Value = 132.870179
$ws.Range("O12").Value = 0.7070964373190639
$ws.Range("P12").Value = 0.7070964373190639
$ws.Range("Q12").Value = 1300.924983812986
$ws.Range("R12").Value = 11708.32485431687
$ws.Range("S12").Value = 0.2948656019495348
$ws.Range("T12").Value = 0.2948656019495349

$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Hgf"
$ws.Range("C13").Value = "Sdc2"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 29.37284333333333
$ws.Range("H13").Value = 88.11852999999999
$ws.Range("I13").Value = 0.417009033544999
$ws.Range("J13").Value = 0.4170090335449991
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 16.09762433333333
$ws.Range("N13").Value = 48.292873
$ws.Range("O13").Value = 0.2570006204793478
$ws.Range("P13").Value = 0.2570006204793479
$ws.Range("Q13").Value = 472.8329975818544
$ws.Range("R13").Value = 4255.49697823669
$ws.Range("S13").Value = 0.1071715803665579
$ws.Range("T13").Value = 0.1071715803665579

$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Hgf"
$ws.Range("C14").Value = "Sdc2"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.2934873333333334
$ws.Range("H14").Value = 0.8804620000000001
$ws.Range("I14").Value = 0.004166667415957767
$ws.Range("J14").Value = 0.004166667415957768
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 2.248835333333334
$ws.Range("N14").Value = 6.746506
$ws.Range("O14").Value = 0.03590294220158827
$ws.Range("P14").Value = 0.03590294220158827
$ws.Range("Q14").Value = 0.6600046850857779
$ws.Range("R14").Value = 5.940042165772001
$ws.Range("S14").Value = 0.0001495956194083729
$ws.Range("T14").Value = 0.0001495956194083729

$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Hgf"
$ws.Range("C15").Value = "Sdc2"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.2934873333333334
$ws.Range("H15").Value = 0.8804620000000001
$ws.Range("I15").Value = 0.004166667415957767
$ws.Range("J15").Value = 0.004166667415957768
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 44.29005966666667
$ws.Range("N15").Value = 132.870179
$ws.Range("O15").Value = 0.7070964373190639
$ws.Range("P15").Value = 0.7070964373190639
$ws.Range("Q15").Value = 12.99857150474423
$ws.Range("R15").Value = 116.987143542698
$ws.Range("S15").Value = 0.002946235685317167
$ws.Range("T15").Value = 0.002946235685317168

$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Hgf"
$ws.Range("C16").Value = "Sdc2"
$ws.Range("D16").Value = "sCs"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.2934873333333334
$ws.Range("H16").Value = 0.8804620000000001
$ws.Range("I16").Value = 0.004166667415957767
$ws.Range("J16").Value = 0.004166667415957768
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 16.09762433333333
$ws.Range("N16").Value = 48.292873
$ws.Range("O16").Value = 0.2570006204793478
$ws.Range("P16").Value = 0.2570006204793479
$ws.Range("Q16").Value = 4.724448838591778
$ws.Range("R16").Value = 42.52003954732601
$ws.Range("S16").Value = 0.001070836111232227
$ws.Range("T16").Value = 0.001070836111232227

